$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 previously held the GlobalBlue URL directly in C5 with no label in D5.
# The edit adds a "GlobalBlue" label in C5 and moves the URL into D5.
$ws.Range("D5").Value = $ws.Range("C5").Value2
$ws.Range("C5").Value = "GlobalBlue"
# Match the vertical-centered formatting used by the other "value" cells in column D.
$ws.Range("D5").VerticalAlignment = -4108

# Widen columns B and D slightly.
$ws.Columns.Item(2).ColumnWidth = 19.2
$ws.Columns.Item(4).ColumnWidth = 110.5

# Move the active selection to G5.
$ws.Range("G5").Select()
